$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Query 1_2 (was using the retired "description = 'Data Analyst'" clause) ---
# Select clause for the first block of queries used to read the old "description = ..." clause;
# it is replaced by the correctly-paired where clause, and the new query now highlights it like
# the other "corrected" cells (C10, C32) with the yellow fill used elsewhere in the sheet.
$ws.Range("C2").Value = 'title = ''Data Analyst'''
$ws.Range("C2").Interior().Color = 65535

# --- Row 3 / Row 4: shift the where-clauses for 1_2 / 1_3 down the list ---
$ws.Range("C3").Value = 'description IN (''Data Analyst'')'
$ws.Range("C4").Value = 'description LIKE ''%Data Analyst%'''

# --- Row 5 (query 1_4): previously missing Select/Where clause, now filled in ---
$ws.Range("B5").Value = '*'
$ws.Range("C5").Value = 'description REGEXP ''Data Analyst'''

# --- Row 10: Query 3_1 (same correction as row 2) ---
$ws.Range("C10").Value = 'title = ''Data Analyst'''
$ws.Range("C10").Interior().Color = 65535

# --- Row 11 / Row 12: shift the where-clauses for 3_2 / 3_3 down the list ---
$ws.Range("C11").Value = 'description IN (''Data Analyst'')'
$ws.Range("C12").Value = 'description LIKE ''%Data Analyst%'''

# --- Row 13 (query 3_4): previously missing Select/Where clause, now filled in ---
$ws.Range("B13").Value = '*'
$ws.Range("C13").Value = 'description REGEXP ''Data Analyst'''

# --- Rows 36-55 (M1-M20): fill in the MongoDB "Where Clause" equivalents ---
$ws.Range("C36").Value = 'find({''description'' : ''Data Analyst''})'
$ws.Range("C37").Value = 'find({''description'' : { ''$in'': [''Data Analyst'']}})'
$ws.Range("C38").Value = 'find({''description'': {''$regex'': ''Data Analyst''}})'
$ws.Range("C39").Value = 'find({''title'' : ''Data Analyst''})'
$ws.Range("C40").Value = 'find({''title'' : { ''$in'': [''Data Analyst'']}})'
$ws.Range("C41").Value = 'find({''title'' : {''$regex'': ''Data Analyst''}})'
$ws.Range("C42").Value = 'find({''title'' : {''$not'' : {''$eq'':''Data Analyst''}}})'
$ws.Range("C43").Value = 'find({''title'' : {''$not'' : {''$regex'': ''Data Analyst''}}})'
$ws.Range("C44").Value = 'find({''$or'' : [{''title'' : ''Data Analyst''} , {''description'' : ''Data Analyst''}]})'
$ws.Range("C45").Value = 'find({''$or'' : [{''title'' : {''$regex'': ''Data Analyst''}} , {''description'' : {''$regex'': ''Data Analyst''}}]})'
$ws.Range("C46").Value = 'find({''$and'' : [{''title'' : ''Data Analyst''} , {''description'' : ''Data Analyst''}]})'
$ws.Range("C47").Value = 'find({''$and'' : [{''title'' : {''$regex'': ''Data Analyst''}} , {''description'' : {''$regex'': ''Data Analyst''}}]})'
$ws.Range("C48").Value = 'find({''title'' : {''$regex'': ''Data Analyst|software engineer''}})'
$ws.Range("C49").Value = 'find({''title'' : {''$regex'': ''Data Analyst.*software engineer''}})'
$ws.Range("C50").Value = 'find({''$or'' : [{''title'' : {''$regex'': ''Data Analyst''}} , {''title'' : {''$regex'': ''software engineer''}}]})'
$ws.Range("C51").Value = 'find({''$and'' : [{''title'' : {''$regex'': ''Data Analyst''}} , {''title'' : {''$regex'': ''software engineer''}}]})'
$ws.Range("C52").Value = 'find({''title'' : ''Data Analyst''},{''title'':1, ''company_name'':1, ''location'':1 })'
$ws.Range("C53").Value = 'find({''title'' : {''$regex'': ''Data Analyst''}},{''title'':1, ''company_name'':1, ''location'':1 })'
$ws.Range("C54").Value = 'find({''description'': {''$regex'': ''Job [a-zA-Z ]*[^a-zA-Z][a-zA-Z ]*Data Analyst''}})'
$ws.Range("C55").Value = 'find({''description'': {''$regex'': ''(Job [a-zA-Z ]*[^a-zA-Z][a-zA-Z ]*(d|D)ata (a|A)nalyst.*(j|J)ob (c|C)lassification[^a-zA-Z].*full.time|(c|C)ompensation:[a-zA-Z ]*(N|n)on[^a-zA-Z ][a-zA-Z ]*$\b(1[0-9])\b[^a-zA-Z ]*([3-9][0-9]))''}})'

# --- Note (F2): replaces the old "fill in clauses" reminder now that they are filled in ---
$ws.Range("F2").Value = 'Clause columns from visualization_analysis file posted 12/1.'

# --- Header cells A1/D1 lose the "Good" (green) highlight style ---
$ws.Range("A1").Style = "Normal"
$ws.Range("D1").Style = "Normal"

# --- Update the active selection to match the author's last cursor position ---
$ws.Range("C6").Select() | Out-Null
